$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content updates -------------------------------------------------
# "Application must behave differently based on logged in user" -> add " role"
$ws.Range("A4").Value = "Application must behave differently based on logged in user role"

# "Application must have an admin page" -> "Application must have an inventory management page"
$ws.Range("A5").Value = "Application must have an inventory management page"

# --- Task status re-shuffle (account creator / admin/inventory page work) -
# Row 5 ("...inventory management page"): move its "X" marker from the
# "In Progress" column (C) to the "To Do" column (B).
$ws.Range("B5").Value = $ws.Range("C5").Value()
$ws.Range("C5").Clear()

# Row 6 ("...order management page"): move its "X" marker from
# "In Progress" (C) to "Completed" (D).
$ws.Range("D6").Value = $ws.Range("C6").Value()
$ws.Range("C6").Clear()

# Row 7 ("...account creation page"): move its "X" marker from
# "To Do" (B) to "Completed" (D).
$ws.Range("D7").Value = $ws.Range("B7").Value()
$ws.Range("B7").ClearContents()

# --- Formatting: the underline emphasis used for "in progress" items is
# dropped now that rows 5 & 6 are no longer "in progress".
$ws.Range("A5").Font.Underline = -4142
$ws.Range("A6").Font.Underline = -4142

# Row 7's label also loses its old "to do" look and adopts the plain look
# used by completed tasks (same as A3 / A9:A12).
$ws.Range("A7").Font.Underline = -4142

# --- Selection -------------------------------------------------------------
$null = $ws.Range("B7").Select()
